$d = $word.ActiveDocument
Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
Write-Output ("Tables: " + $d.Tables.Count)
